$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.623.34"
$ws.Range("E2").Value = "  +3.76%  "

$ws.Range("D3").Value = "1.919.66"
$ws.Range("E3").Value = "  +2.16%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.06"
$ws.Range("E5").Value = "  +1.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.700"
$ws.Range("E6").Value = "  +2.56%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.51"
$ws.Range("E8").Value = "  +3.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "58.88"
$ws.Range("E9").Value = "  +10.20%  "

$ws.Range("E10").Value = "  +3.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0764"
$ws.Range("E11").Value = "  +3.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.100"
$ws.Range("E12").Value = "  +3.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.58"
$ws.Range("E13").Value = "  +8.42%  "

$ws.Range("E14").Value = "  +3.87%  "

$ws.Range("D15").Value = "2.198.90"
$ws.Range("E15").Value = "  +2.17%  "

$ws.Range("E16").Value = "  +4.80%  "

$ws.Range("D17").Value = "1.915.62"
$ws.Range("E17").Value = "  +2.06%  "

$ws.Range("D18").Value = "36.685.87"
$ws.Range("E18").Value = "  +3.86%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.16"
$ws.Range("E19").Value = "  +2.03%  "

$ws.Range("D20").Value = "0.0₃0861"
$ws.Range("E20").Value = "  +5.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "251.68"
$ws.Range("E21").Value = "  +3.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.26"
$ws.Range("E22").Value = "  +4.08%  "

$ws.Range("E23").Value = "  +5.23%  "

$ws.Range("E24").Value = "  +1.76%  "

$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("E26").Value = "  +2.32%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.88"
$ws.Range("E27").Value = "  +1.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.82"
$ws.Range("E28").Value = "  +3.68%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.77"
$ws.Range("E29").Value = "  +3.03%  "

$ws.Range("E30").Value = "  +2.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.54"
$ws.Range("E31").Value = "  +6.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0625"
$ws.Range("E32").Value = "  +6.79%  "

$ws.Range("E33").Value = "  +0.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.36"
$ws.Range("E34").Value = "  +5.73%  "

$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0848"
$ws.Range("E36").Value = "  +15.74%  "

$ws.Range("E37").Value = "  -12.47%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.874"
$ws.Range("E38").Value = "  +4.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.90"
$ws.Range("E39").Value = "  +48.59%  "

$ws.Range("E40").Value = "  +4.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.15"
$ws.Range("E41").Value = "  +10.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0228"
$ws.Range("E42").Value = "  +5.83%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.24"
$ws.Range("E43").Value = "  -2.20%  "

$ws.Range("E44").Value = "  +3.39%  "

$ws.Range("D45").Value = "1.337.56"
$ws.Range("E45").Value = "  +2.71%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.55"
$ws.Range("E46").Value = "  +7.05%  "

$ws.Range("E47").Value = "  +1.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0815"
$ws.Range("E48").Value = "  +2.48%  "

$ws.Range("E49").Value = "  +2.49%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.46"
$ws.Range("E50").Value = "  +4.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.10"
$ws.Range("E51").Value = "  +2.69%  "
